$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# The existing row 11's created_at timestamp is refreshed to a (very
# slightly) later value, as if re-written by the application.
$ws.Range("I11").Value = 45790.74601511574

# A new user registration row is appended to the Users table.
$row = 12
$ws.Cells.Item($row, 1).Value = "c5578616-d725-4cc3-a2f5-4f5b5ce78ac2"
$ws.Cells.Item($row, 2).Value = "mayy"
$ws.Cells.Item($row, 3).Value = "gerald.mandebvu@gmail.com"
$ws.Cells.Item($row, 4).Value = ""
$ws.Cells.Item($row, 5).Value = "scrypt:32768:8:1`$orxKA4edlVvQvpyn`$83571dfbe4eee3d19890937bfbacac0b1c7a438aa4607398e299ae711d706781d2a457ac717f1e85a73585601b4a1dc5c072e812fa05b2f4f10184ae4bbc6e70"
$ws.Cells.Item($row, 6).Value = "ITRADE-25480410"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = ""
$ws.Cells.Item($row, 9).Value = 45790.86667068065
$ws.Cells.Item($row, 9).NumberFormat = $ws.Cells.Item(11, 9).NumberFormat
